$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new columns to make room for new HPO-term / data columns ---
# (inserted left-to-right so each Insert() targets the correct final position)
$ws.Columns("T:U").Insert()   # new: Atopic dermatitis, Erythematous plaque
$ws.Columns("W:W").Insert()   # new: Trichorrhexis invaginata
$ws.Columns("Y:Y").Insert()   # new: Dermatographic urticaria
$ws.Columns("AN:AR").Insert() # new: Abnormal IgG/IgA/IgM concentration, Elevated tryptase, Eosinophilia

# --- Row 1: headers for newly inserted columns ---
$ws.Range("T1").Value = 'Atopic dermatitis'
$ws.Range("U1").Value = 'Erythematous plaque'
$ws.Range("W1").Value = 'Trichorrhexis invaginata'
$ws.Range("Y1").Value = 'Dermatographic urticaria'
$ws.Range("AN1").Value = 'Abnormal circulating IgG concentration'
$ws.Range("AO1").Value = 'Abnormal circulating IgA concentration'
$ws.Range("AP1").Value = 'Abnormal circulating IgM concentration'
$ws.Range("AQ1").Value = 'Elevated total serum tryptase'
$ws.Range("AR1").Value = 'Eosinophilia'

# --- Row 1: new trailing columns appended after BA (Eclabion) ---
$ws.Range("BB1").Value = 'Nut food product allergy'
$ws.Range("BC1").Value = 'Nasal congestion'

# --- Row 2: CURIE / type row for newly inserted + appended columns ---
$ws.Range("T2").Value = 'HP:0001047'
$ws.Range("U2").Value = 'HP:0025474'
$ws.Range("W2").Value = 'HP:0025811'
$ws.Range("Y2").Value = 'HP:0011971'
$ws.Range("AN2").Value = 'HP:0410242'
$ws.Range("AO2").Value = 'HP:0410240'
$ws.Range("AP2").Value = 'HP:0410243'
$ws.Range("AQ2").Value = 'HP:0031901'
$ws.Range("AR2").Value = 'HP:0001880'
$ws.Range("BB2").Value = 'HP:0410331'
$ws.Range("BC2").Value = 'HP:0001742'

# --- Row 4: new patient record (PMID:39891497 - SPINK5 variants case report) ---
$ws.Range("A4").Value = 'PMID:39891497'
$ws.Range("B4").Value = 'Systemic JAK inhibitors for treatment of cutaneous manifestations in a patient with SPINK5 variants: A case report and review of the literature'
$ws.Range("C4").Value = 'patient'
$ws.Range("E4").Value = 'OMIM:256500'
$ws.Range("F4").Value = 'Netherton syndrome'
$ws.Range("G4").Value = 'HGNC:15464'
$ws.Range("H4").Value = 'SPINK5'
$ws.Range("I4").Value = 'NM_006846.4'
$ws.Range("J4").Value = 'c.2390G>T'
$ws.Range("K4").Value = 'c.1499G>A'
$ws.Range("L4").Value = 'NP_006837.2:p.(Gly797Val);NP_006837.2:p.(Arg500Gln)'
$ws.Range("M4").Value = 'P1M'
$ws.Range("N4").Value = 'P17Y'
$ws.Range("O4").Value = 'no'
$ws.Range("P4").Value = 'M'
$ws.Range("Q4").Value = 'na'
$ws.Range("T4").Value = 'P1M'
$ws.Range("U4").Value = 'observed'
$ws.Range("W4").Value = 'observed'
$ws.Range("Y4").Value = 'observed'
$ws.Range("AM4").Value = 'observed'
$ws.Range("AN4").Value = 'excluded'
$ws.Range("AO4").Value = 'excluded'
$ws.Range("AP4").Value = 'excluded'
$ws.Range("AQ4").Value = 'excluded'
$ws.Range("AR4").Value = 'excluded'
$ws.Range("BB4").Value = 'observed'
$ws.Range("BC4").Value = 'observed'
